$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5558125
$ws.Range("I132").Value = 6668603.5
$ws.Range("J132").Value = 5733.3335
$ws.Range("K132").Value = 20005810.5
$ws.Range("L132").Value = 17200.0005
$ws.Range("M132").Value = -20003280.5
$ws.Range("N132").Value = -22260.0005
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").Value = 0
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("N134").Value = 0
$ws.Range("H137").Value = 1889.6
$ws.Range("I137").Value = 1377.2727
$ws.Range("J137").Value = 2186.2104
$ws.Range("K137").Value = 4131.8181
$ws.Range("L137").Value = 6558.6312
$ws.Range("M137").Value = -1581.8181
$ws.Range("N137").Value = -11658.6312
$ws.Range("H139").Value = 30000
$ws.Range("J139").Value = 30000
$ws.Range("L139").Value = 30000
$ws.Range("N139").Value = -40280
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").ClearContents()
$ws.Range("N140").Value = 0
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5642.719
$ws.Range("I32").Value = 4137.814
$ws.Range("J32").Value = 10264.929
$ws.Range("K32").Value = 4137.814
$ws.Range("L32").Value = 10264.929
$ws.Range("M32").Value = -3850.814
$ws.Range("N32").Value = -10838.929
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").ClearContents()
$ws.Range("N44").Value = 0
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H61").Value = 2716.1191
$ws.Range("I61").Value = 951.05
$ws.Range("K61").Value = 951.05
$ws.Range("M61").Value = -739.05
$ws.Range("H80").Value = 31205
$ws.Range("J80").Value = 31205
$ws.Range("L80").Value = 31205
$ws.Range("N80").Value = -33201
$ws.Range("H83").Value = 31205
$ws.Range("J83").Value = 31205
$ws.Range("L83").Value = 93615
$ws.Range("N83").Value = -103599
$ws.Range("H132").Value = 45460420
$ws.Range("I132").Value = 76929896
$ws.Range("J132").Value = 4510.8887
$ws.Range("K132").Value = 230789688
$ws.Range("L132").Value = 13532.6661
$ws.Range("M132").Value = -230787158
$ws.Range("N132").Value = -18592.6661
$ws.Range("H136").Value = 2716.1191
$ws.Range("I136").Value = 951.05
$ws.Range("K136").Value = 2853.15
$ws.Range("M136").Value = -303.1499999999996
$ws.Range("H137").Value = 39333.332
$ws.Range("J137").Value = 39333.332
$ws.Range("L137").Value = 39333.332
$ws.Range("N137").Value = -49533.332
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H119").Value = 25000
$ws.Range("J119").Value = 25000
$ws.Range("L119").Value = 25000
$ws.Range("N119").Value = -34676
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 32400.8
$ws.Range("J4").Value = 50668
$ws.Range("L4").Value = 50668
$ws.Range("N4").Value = -50892
$ws.Range("H31").Value = 2227.7463
$ws.Range("I31").Value = 1218.5227
$ws.Range("J31").Value = 4158.4346
$ws.Range("K31").Value = 1218.5227
$ws.Range("L31").Value = 4158.4346
$ws.Range("M31").Value = -923.5227
$ws.Range("N31").Value = -4748.4346
$ws.Range("H34").Value = 2227.7463
$ws.Range("I34").Value = 1218.5227
$ws.Range("J34").Value = 4158.4346
$ws.Range("K34").Value = 1218.5227
$ws.Range("L34").Value = 4158.4346
$ws.Range("M34").Value = -1016.5227
$ws.Range("N34").Value = -4562.4346
$ws.Range("H134").Value = 1747.2881
$ws.Range("I134").Value = 1362.9375
$ws.Range("J134").Value = 3424.4546
$ws.Range("K134").Value = 4088.8125
$ws.Range("L134").Value = 10273.3638
$ws.Range("M134").Value = -1553.8125
$ws.Range("N134").Value = -15343.3638
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 131.71428
$ws.Range("I2").Value = 150.75
$ws.Range("J2").Value = 124.1
$ws.Range("K2").Value = 904.5
$ws.Range("L2").Value = 744.5999999999999
$ws.Range("M2").Value = -791.5
$ws.Range("N2").Value = -970.5999999999999
$ws.Range("H14").Value = 331.33334
$ws.Range("I14").Value = 331.33334
$ws.Range("K14").Value = 994.0000200000001
$ws.Range("M14").Value = -821.0000200000001
$ws.Range("H23").Value = 158
$ws.Range("I23").Value = 95
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 285
$ws.Range("L23").Value = 600
$ws.Range("M23").Value = -50
$ws.Range("N23").Value = -1070
$ws.Range("H34").Value = 15642.857
$ws.Range("I34").Value = 450
$ws.Range("J34").Value = 21720
$ws.Range("K34").Value = 1350
$ws.Range("L34").Value = 65160
$ws.Range("M34").Value = -1266
$ws.Range("N34").Value = -65328
$ws.Range("H39").Value = 3450
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 3450
$ws.Range("K39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("M39").Value = 10350
$ws.Range("N39").Value = -10938
$ws.Range("H87").Value = 12966.667
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 12966.667
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H107").Value = 967.54236
$ws.Range("I107").Value = 650.41174
$ws.Range("J107").Value = 1398.84
$ws.Range("K107").Value = 1951.23522
$ws.Range("L107").Value = 4196.52
$ws.Range("M107").Value = -31.23522000000003
$ws.Range("N107").Value = -8036.52
$ws.Range("H118").Value = 4078.2173
$ws.Range("I118").Value = 724.75
$ws.Range("K118").Value = 2174.25
$ws.Range("M118").Value = -931.25
$ws.Range("H126").Value = 2076.6667
$ws.Range("I126").Value = 1230
$ws.Range("K126").Value = 3690
$ws.Range("M126").Value = 1250
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 52687.332
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 52687.332
$ws.Range("K44").Value = 0
$ws.Range("L44").ClearContents()
$ws.Range("M44").Value = 52687.332
$ws.Range("N44").Value = -53879.332
$ws.Range("H126").Value = 3121.4285
$ws.Range("I126").Value = 1966.6666
$ws.Range("J126").Value = 4453.846
$ws.Range("K126").Value = 5899.9998
$ws.Range("L126").Value = 13361.538
$ws.Range("M126").Value = -3429.9998
$ws.Range("N126").Value = -18301.538
$ws.Range("H132").Value = 3530.423
$ws.Range("I132").Value = 2469
$ws.Range("J132").Value = 4768.75
$ws.Range("K132").Value = 7407
$ws.Range("L132").Value = 14306.25
$ws.Range("M132").Value = -4877
$ws.Range("N132").Value = -19366.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 166670370
$ws.Range("I22").Value = 333333980
$ws.Range("J22").Value = 6733.3335
$ws.Range("K22").Value = 333333980
$ws.Range("L22").Value = 6733.3335
$ws.Range("M22").Value = -333333685
$ws.Range("N22").Value = -7323.3335
$ws.Range("H27").Value = 166670370
$ws.Range("I27").Value = 333333980
$ws.Range("J27").Value = 6733.3335
$ws.Range("K27").Value = 333333980
$ws.Range("L27").Value = 6733.3335
$ws.Range("M27").Value = -333333873
$ws.Range("N27").Value = -6947.3335
$ws.Range("H55").Value = 988.7778
$ws.Range("I55").Value = 415.6
$ws.Range("J55").Value = 1209.2307
$ws.Range("K55").Value = 415.6
$ws.Range("L55").Value = 1209.2307
$ws.Range("M55").Value = -242.6
$ws.Range("N55").Value = -1555.2307
$ws.Range("H132").Value = 2041.8462
$ws.Range("I132").Value = 1315.25
$ws.Range("J132").Value = 3891.3635
$ws.Range("K132").Value = 3945.75
$ws.Range("L132").Value = 11674.0905
$ws.Range("M132").Value = -1415.75
$ws.Range("N132").Value = -16734.0905
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4349604
$ws.Range("I126").Value = 1183.4615
$ws.Range("J126").Value = 10002551
$ws.Range("K126").Value = 3550.3845
$ws.Range("L126").Value = 30007653
$ws.Range("M126").Value = -1080.3845
$ws.Range("N126").Value = -30012593
$ws.Range("H132").Value = 3387.5
$ws.Range("I132").Value = 1169.3556
$ws.Range("J132").Value = 9259.058999999999
$ws.Range("K132").Value = 3508.066800000001
$ws.Range("L132").Value = 27777.177
$ws.Range("M132").Value = -978.0668000000005
$ws.Range("N132").Value = -32837.177
